$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = [string]([char]0x00C2) + [string]([char]0x00B1)
$replacement = [string]([char]0x00B1)

for ($r = 2; $r -le 17; $r++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value()
        if ($val -ne $null -and $val.Contains($target)) {
            $newVal = $val.Replace($target, $replacement)
            $cell.Value = $newVal
        }
    }
}
